# TC08_C3DC_phs002517_DiseasePhase-NotReported.xlsx
# "Updated remaining queries for C3DC"
#
# The sheet holds a set of SQL queries (one per cell) that join several
# DuckDB/pandas-backed views together. Previously every join used the
# generic "id" columns (std.id / prt.id) together with dotted aliases
# ("study.id" / "participant.id"). The queries are updated so the joins
# use the fully qualified id columns instead (std.study_id / prt.participant_id,
# "study.study_id" / "participant.participant_id").
#
# This touches every query cell on the sheet: C2 (StatQuery/count query),
# B2 (Studies), B3 (Participants), B4 (Diagnosis), B5 (Treatment),
# B6 (Treatment Response) and B7 (Survival).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Update-QueryJoins([string]$text) {
    $text = $text.Replace('std.id = prt."study.id"', 'std.study_id = prt."study.study_id"')
    $text = $text.Replace('prt.id = dgn."participant.id"', 'prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('prt.id = trt."participant.id"', 'prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('prt.id = trr."participant.id"', 'prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('prt.id = srv."participant.id"', 'prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('std.id = rfs."study.id"', 'std.study_id = rfs."study.study_id"')
    return $text
}

$cellsToFix = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")
foreach ($addr in $cellsToFix) {
    $range = $ws.Range($addr)
    $current = $range.Value()
    $updated = Update-QueryJoins $current
    $range.Value = $updated
}

# Column C (StatQuery) was manually widened while reviewing the updated
# queries, which also clears the previous "best fit" auto-sizing.
$ws.Columns.Item(3).ColumnWidth = 67.33

# The workbook was left scrolled down with C7 (Survival query cell) selected.
$ws.Activate()
$ws.Range("C7").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
